$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 117.5
$ws.Range("I5").Value = 27.5
$ws.Range("J5").Value = 162.5
$ws.Range("K5").Value = 27.5
$ws.Range("L5").Value = 162.5
$ws.Range("M5").Value = 87.5
$ws.Range("N5").Value = -392.5

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = $null

$ws.Range("H137").Value = 6745.857
$ws.Range("I137").Value = 9963.666999999999
$ws.Range("J137").Value = 2455.4443
$ws.Range("K137").Value = 29891.001
$ws.Range("L137").Value = 7366.3329
$ws.Range("M137").Value = -27341.001
$ws.Range("N137").Value = -12466.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I4").Value = 133.33333
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 133.33333
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -17.33332999999999
$ws.Range("N4").Value = -432

$ws.Range("H5").Value = 145.5
$ws.Range("I5").Value = 93.333336
$ws.Range("J5").Value = 302
$ws.Range("K5").Value = 93.333336
$ws.Range("L5").Value = 302
$ws.Range("M5").Value = 18.666664
$ws.Range("N5").Value = -526

$ws.Range("H32").Value = 7606.3335
$ws.Range("I32").Value = 3738.0417
$ws.Range("J32").Value = 19984.867
$ws.Range("K32").Value = 3738.0417
$ws.Range("L32").Value = 19984.867
$ws.Range("M32").Value = -3451.0417
$ws.Range("N32").Value = -20558.867

$ws.Range("H61").Value = 1958.0193
$ws.Range("I61").Value = 1684.6
$ws.Range("J61").Value = 2869.4167
$ws.Range("K61").Value = 1684.6
$ws.Range("L61").Value = 2869.4167
$ws.Range("M61").Value = -1472.6
$ws.Range("N61").Value = -3293.4167

$ws.Range("H74").Value = 4551226
$ws.Range("I74").Value = 9091640
$ws.Range("J74").Value = 10811.909
$ws.Range("K74").Value = 9091640
$ws.Range("L74").Value = 10811.909
$ws.Range("M74").Value = -9090766
$ws.Range("N74").Value = -12559.909

$ws.Range("H77").Value = 4551226
$ws.Range("I77").Value = 9091640
$ws.Range("J77").Value = 10811.909
$ws.Range("K77").Value = 45458200
$ws.Range("L77").Value = 54059.545
$ws.Range("M77").Value = -45453832
$ws.Range("N77").Value = -62795.545

$ws.Range("H136").Value = 1958.0193
$ws.Range("I136").Value = 1684.6
$ws.Range("J136").Value = 2869.4167
$ws.Range("K136").Value = 5053.799999999999
$ws.Range("L136").Value = 8608.250100000001
$ws.Range("M136").Value = -2503.799999999999
$ws.Range("N136").Value = -13708.2501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 145.5
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 302
$ws.Range("K4").Value = 93.333336
$ws.Range("L4").Value = 302
$ws.Range("M4").Value = 21.666664
$ws.Range("N4").Value = -532

$ws.Range("H107").Value = 2742.0527
$ws.Range("I107").Value = 2117.3
$ws.Range("J107").Value = 3436.2222
$ws.Range("K107").Value = 2117.3
$ws.Range("L107").Value = 3436.2222
$ws.Range("M107").Value = -197.3000000000002
$ws.Range("N107").Value = -7276.2222

$ws.Range("H134").Value = 31062.918
$ws.Range("I134").Value = 34498.426
$ws.Range("J134").Value = 2720
$ws.Range("K134").Value = 103495.278
$ws.Range("L134").Value = 8160
$ws.Range("M134").Value = -100960.278
$ws.Range("N134").Value = -13230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2114.8
$ws.Range("I7").Value = 2021.1428
$ws.Range("J7").Value = 2333.3333
$ws.Range("K7").Value = 2021.1428
$ws.Range("L7").Value = 2333.3333
$ws.Range("M7").Value = -1908.1428
$ws.Range("N7").Value = -2559.3333

$ws.Range("H31").Value = 1564.1428
$ws.Range("I31").Value = 1010.9091
$ws.Range("J31").Value = 2500.3845
$ws.Range("K31").Value = 1010.9091
$ws.Range("L31").Value = 2500.3845
$ws.Range("M31").Value = -715.9091
$ws.Range("N31").Value = -3090.3845

$ws.Range("H34").Value = 1564.1428
$ws.Range("I34").Value = 1010.9091
$ws.Range("J34").Value = 2500.3845
$ws.Range("K34").Value = 1010.9091
$ws.Range("L34").Value = 2500.3845
$ws.Range("M34").Value = -808.9091
$ws.Range("N34").Value = -2904.3845

$ws.Range("H122").Value = 31250804
$ws.Range("I122").Value = 31250804
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 93752412
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -93749962

$ws.Range("H132").Value = 4149.4116
$ws.Range("I132").Value = 4003.2307
$ws.Range("J132").Value = 4624.5
$ws.Range("K132").Value = 12009.6921
$ws.Range("L132").Value = 13873.5
$ws.Range("M132").Value = -9479.6921
$ws.Range("N132").Value = -18933.5

$ws.Range("H134").Value = 1930.5555
$ws.Range("I134").Value = 1716.4667
$ws.Range("J134").Value = 3001
$ws.Range("K134").Value = 5149.4001
$ws.Range("L134").Value = 9003
$ws.Range("M134").Value = -2614.4001
$ws.Range("N134").Value = -14073

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 129.125
$ws.Range("I11").Value = 76.14286
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 228.42858
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = -88.42858000000001
$ws.Range("N11").Value = -1780

$ws.Range("H34").Value = 66667200
$ws.Range("I34").Value = 255.55556
$ws.Range("J34").Value = 166667620
$ws.Range("K34").Value = 766.66668
$ws.Range("L34").Value = 500002860
$ws.Range("M34").Value = -682.66668
$ws.Range("N34").Value = -500003028

$ws.Range("H131").Value = 1540362.6
$ws.Range("I131").Value = 2286
$ws.Range("J131").Value = 1668535.8
$ws.Range("K131").Value = 6858
$ws.Range("L131").Value = 5005607.4
$ws.Range("M131").Value = -1818
$ws.Range("N131").Value = -5015687.4

$ws.Range("H132").Value = 62501040
$ws.Range("I132").Value = 250001700
$ws.Range("J132").Value = 818.6667
$ws.Range("K132").Value = 2250015300
$ws.Range("L132").Value = 7368.0003
$ws.Range("M132").Value = -2250012770
$ws.Range("N132").Value = -12428.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = $null
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = 0

$ws.Range("H22").Value = 530.7273
$ws.Range("I22").Value = 535.8
$ws.Range("J22").Value = 480
$ws.Range("K22").Value = 535.8
$ws.Range("L22").Value = 480
$ws.Range("M22").Value = -240.8
$ws.Range("N22").Value = -1070

$ws.Range("H27").Value = 530.7273
$ws.Range("I27").Value = 535.8
$ws.Range("J27").Value = 480
$ws.Range("K27").Value = 535.8
$ws.Range("L27").Value = 480
$ws.Range("M27").Value = -428.8
$ws.Range("N27").Value = -694

$ws.Range("H34").Value = 12000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 12000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 10775.25
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 10775.25
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 10775.25
$ws.Range("N12").Value = -11059.25

$ws.Range("H40").Value = 10998.667
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 14998
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 14998
$ws.Range("M40").Value = -2851
$ws.Range("N40").Value = -15296
